$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Inscricoes")

# Row 5: Inscritos 104 -> 106
$ws.Range("E5").Value = 106

# Row 10: Inscritos 355 -> 357
$ws.Range("E10").Value = 357

# Row 12: Inscritos 350 -> 351
$ws.Range("E12").Value = 351

# Row 15: Inscritos 118 -> 119
$ws.Range("E15").Value = 119

# Row 16: Inscritos 156 -> 157
$ws.Range("E16").Value = 157

# Row 24: Inscritos 156 -> 158, Pagos 76 -> 77, Inscricoes homologadas 76 -> 77
$ws.Range("E24").Value = 158
$ws.Range("F24").Value = 77
$ws.Range("H24").Value = 77

# Row 26: Inscritos 99 -> 100, Pagos 59 -> 60, Inscricoes homologadas 59 -> 60
$ws.Range("E26").Value = 100
$ws.Range("F26").Value = 60
$ws.Range("H26").Value = 60

# Row 27: Inscritos 239 -> 241
$ws.Range("E27").Value = 241

# Row 32: Inscritos 142 -> 143, Pagos 76 -> 77, Inscricoes homologadas 76 -> 77
$ws.Range("E32").Value = 143
$ws.Range("F32").Value = 77
$ws.Range("H32").Value = 77

# Row 34: Inscritos 162 -> 165, Pagos 95 -> 96, Inscricoes homologadas 95 -> 96
$ws.Range("E34").Value = 165
$ws.Range("F34").Value = 96
$ws.Range("H34").Value = 96

# Row 41: Inscritos 297 -> 299
$ws.Range("E41").Value = 299

# Row 42: Inscritos 263 -> 264, Pagos 137 -> 139, Inscricoes homologadas 137 -> 139
$ws.Range("E42").Value = 264
$ws.Range("F42").Value = 139
$ws.Range("H42").Value = 139

# Row 43: Inscritos 87 -> 89, Pagos 41 -> 42, Inscricoes homologadas 41 -> 42
$ws.Range("E43").Value = 89
$ws.Range("F43").Value = 42
$ws.Range("H43").Value = 42

# Row 44: Inscritos 243 -> 244
$ws.Range("E44").Value = 244

# Row 45: Inscritos 101 -> 102
$ws.Range("E45").Value = 102

# Row 46: Inscritos 230 -> 231
$ws.Range("E46").Value = 231

# Row 50: Inscritos 189 -> 190, Pagos 67 -> 68, Inscricoes homologadas 67 -> 68
$ws.Range("E50").Value = 190
$ws.Range("F50").Value = 68
$ws.Range("H50").Value = 68

# Row 51: Inscritos 180 -> 181
$ws.Range("E51").Value = 181
